# Restore revision: update cell C10 on the "Rules" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
if ($ws -eq $null) { $ws = $wb.ActiveSheet }

$ws.Range("C10").Value = 1
